# Applies the "Adding exception to create profile, creating UI managment profile"
# update to the Tablut work journal worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Tweak existing comments -------------------------------------------------

# Row 20 (E20): add the missing trailing period.
$ws.Range("E20").Value = "La partie analyse n'est pas encore terminée."

# Row 24 (E24): note the duplicate-search exception, and give it the same
# "multi-line comment" styling already used on row 14 (wrap text over a taller row).
$ws.Range("E14:G14").Copy($ws.Range("E24:G24")) | Out-Null
$ws.Range("E24").Value = "Tout a été implémenté correctement sauf la recherche de doublon."

# --- 2. New row 25: creation of profile (with exception handling) ---------------

# Row 25 already has the right "Task" styling (C25/D25 unchanged); only the date
# cell and the (now taller, wrap-text) comment cell need a style change, so pull
# those specifically from row 14, which already uses the desired formatting.
$ws.Range("B14").Copy($ws.Range("B25")) | Out-Null
$ws.Range("E14:G14").Copy($ws.Range("E25:G25")) | Out-Null
$ws.Rows.Item(25).RowHeight = 47.25

$ws.Range("B25").Value = 43147
$ws.Range("D25").Value = "30 min"
$ws.Range("E25").Value = "Permet de savoir si un profil existe déjà, si la syntax du nom est correct ou si le serveur de base de données est allumé."
$ws.Range("C25").Value = "Création de profil."

# --- 3. New row 26: creation of the profile management UI -----------------------

# Row 26 already carries the right "Comments" styling (E26:G26 blank, bordered);
# only the date cell needs the "date" style used throughout column B.
$ws.Range("B19").Copy($ws.Range("B26")) | Out-Null

$ws.Range("B26").Value = 43147
$ws.Range("C26").Value = "Création UI gestion de profil."
$ws.Range("D26").Value = "1h"

# --- 4. Restore the selection that Excel had when the sheet was last saved ------

$ws.Range("D26").Select() | Out-Null
